$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.470157689643088
$ws.Range("C2").Value = 0.2339133357980074
$ws.Range("D2").Value = 0.1075226812809156
$ws.Range("E2").Value = 0.05439836397513442
$ws.Range("F2").Value = 2.312878378896855
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 1.736421287184115
$ws.Range("L2").Value = 0.2243747179933493
$ws.Range("M2").Value = 0.3159931720354976
$ws.Range("B3").Value = 1.386189995712527
$ws.Range("C3").Value = 0.2037510629328096
$ws.Range("D3").Value = 0.1076720869938619
$ws.Range("E3").Value = 0.05394962065763487
$ws.Range("F3").Value = 2.259330592931931
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 1.714308406271002
$ws.Range("L3").Value = 0.2214141773817815
$ws.Range("M3").Value = 0.3028853873112425
$ws.Range("B4").Value = 1.335603170010813
$ws.Range("C4").Value = 0.1852597149878932
$ws.Range("D4").Value = 0.1077912464955233
$ws.Range("E4").Value = 0.05366683351487023
$ws.Range("F4").Value = 2.227720188075963
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 1.701542869023456
$ws.Range("L4").Value = 0.2197126688251032
$ws.Range("M4").Value = 0.2950430642277553
$ws.Range("B5").Value = 1.31523181244512
$ws.Range("C5").Value = 0.1777308649236033
$ws.Range("D5").Value = 0.1078466504844968
$ws.Range("E5").Value = 0.05354975700494347
$ws.Range("F5").Value = 2.215155243583993
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 1.696543617103472
$ws.Range("L5").Value = 0.2190484816087093
$ws.Range("M5").Value = 0.2918989163574182
$ws.Range("B6").Value = 1.311863833675943
$ws.Range("C6").Value = 0.1764810812518363
$ws.Range("D6").Value = 0.1078562622551402
$ws.Range("E6").Value = 0.05353020510576734
$ws.Range("F6").Value = 2.213087895109012
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 1.695725705140845
$ws.Range("L6").Value = 0.218939955574605
$ws.Range("M6").Value = 0.2913799500532051
$ws.Range("B7").Value = 1.335327450952946
$ws.Range("C7").Value = 0.1851581526719315
$ws.Range("D7").Value = 0.1077919660431519
$ws.Range("E7").Value = 0.05366526204128963
$ws.Range("F7").Value = 2.227549454237163
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 1.701474627992212
$ws.Range("L7").Value = 0.2197035932298519
$ws.Range("M7").Value = 0.2950004521198935
$ws.Range("B8").Value = 1.441003940586995
$ws.Range("C8").Value = 0.2235071036622855
$ws.Range("D8").Value = 0.1075684741528349
$ws.Range("E8").Value = 0.05424513512342877
$ws.Range("F8").Value = 2.294150467946025
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 1.72862744622384
$ws.Range("L8").Value = 0.2233297603923319
$ws.Range("M8").Value = 0.3114307952182003
$ws.Range("B9").Value = 1.655972928619065
$ws.Range("C9").Value = 0.2989665458969739
$ws.Range("D9").Value = 0.1073500978287072
$ws.Range("E9").Value = 0.05532533564858078
$ws.Range("F9").Value = 2.434931384048951
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 1.788380427314422
$ws.Range("L9").Value = 0.2313660915370974
$ws.Range("M9").Value = 0.3452920438299572
$ws.Range("B10").Value = 1.818709298254703
$ws.Range("C10").Value = 0.3546131732741742
$ws.Range("D10").Value = 0.1073269170643982
$ws.Range("E10").Value = 0.05608520854000876
$ws.Range("F10").Value = 2.544738275459338
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 1.836342906360784
$ws.Range("L10").Value = 0.2378395095869905
$ws.Range("M10").Value = 0.3711845248298928
$ws.Range("B11").Value = 1.893803776385482
$ws.Range("C11").Value = 0.3799843339752442
$ws.Range("D11").Value = 0.1073468735177414
$ws.Range("E11").Value = 0.05642377703129053
$ws.Range("F11").Value = 2.596114556266713
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 1.859065365430084
$ws.Range("L11").Value = 0.2409091785256123
$ws.Range("M11").Value = 0.3831873330131543
$ws.Range("B12").Value = 1.922394475568353
$ws.Range("C12").Value = 0.389600772480037
$ws.Range("D12").Value = 0.1073588733534763
$ws.Range("E12").Value = 0.05655097919170871
$ws.Range("F12").Value = 2.615777120314391
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 1.867801379878344
$ws.Range("L12").Value = 0.2420896178971788
$ws.Range("M12").Value = 0.387764936968253
$ws.Range("B13").Value = 1.916230090375279
$ws.Range("C13").Value = 0.387529291901501
$ws.Range("D13").Value = 0.1073560904750295
$ws.Range("E13").Value = 0.05652362845154268
$ws.Range("F13").Value = 2.611533165379569
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 1.865914049695007
$ws.Range("L13").Value = 0.241834586550965
$ws.Range("M13").Value = 0.3867776243353305
$ws.Range("B14").Value = 1.8961528575839
$ws.Range("C14").Value = 0.3807753011134309
$ws.Range("D14").Value = 0.1073477713857116
$ws.Range("E14").Value = 0.05643426211338554
$ws.Range("F14").Value = 2.597728033850359
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 1.859781438224246
$ws.Range("L14").Value = 0.2410059324084415
$ws.Range("M14").Value = 0.3835632849215216
$ws.Range("B15").Value = 1.883875062757397
$ws.Range("C15").Value = 0.3766394707894278
$ws.Range("D15").Value = 0.1073432559704486
$ws.Range("E15").Value = 0.05637939203385578
$ws.Range("F15").Value = 2.589299090884055
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 1.856042208314122
$ws.Range("L15").Value = 0.2405007068276177
$ws.Range("M15").Value = 0.3815986337473092
$ws.Range("B16").Value = 1.813823203404638
$ws.Range("C16").Value = 0.3529563195008905
$ws.Range("D16").Value = 0.1073262314680719
$ws.Range("E16").Value = 0.05606294065320139
$ws.Range("F16").Value = 2.541409604931943
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 1.834876261190274
$ws.Range("L16").Value = 0.2376414178855981
$ws.Range("M16").Value = 0.3704046380943211
$ws.Range("B17").Value = 1.77112204036564
$ws.Range("C17").Value = 0.338442626053677
$ws.Range("D17").Value = 0.1073236398622441
$ws.Range("E17").Value = 0.0558669995827783
$ws.Range("F17").Value = 2.512397686734687
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 1.822124263937155
$ws.Range("L17").Value = 0.2359193666987096
$ws.Range("M17").Value = 0.3635950180075724
$ws.Range("B18").Value = 1.74666163515343
$ws.Range("C18").Value = 0.3301000775852572
$ws.Range("D18").Value = 0.1073250167745172
$ws.Range("E18").Value = 0.05575363013559809
$ws.Range("F18").Value = 2.49584472770843
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 1.814874658575263
$ws.Range("L18").Value = 0.234940638555571
$ws.Range("M18").Value = 0.3596993961436397
$ws.Range("B19").Value = 1.73839693997337
$ws.Range("C19").Value = 0.3272763307396644
$ws.Range("D19").Value = 0.1073259737535679
$ws.Range("E19").Value = 0.05571512978417115
$ws.Range("F19").Value = 2.490263102934847
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 1.81243462096873
$ws.Range("L19").Value = 0.2346112741771407
$ws.Range("M19").Value = 0.3583840231585853
$ws.Range("B20").Value = 1.775657277879759
$ws.Range("C20").Value = 0.3399870753638083
$ws.Range("D20").Value = 0.1073236185544957
$ws.Range("E20").Value = 0.05588792701033363
$ws.Range("F20").Value = 2.515472174854665
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 1.823472927132272
$ws.Range("L20").Value = 0.2361014654278364
$ws.Range("M20").Value = 0.3643177296442204
$ws.Range("B21").Value = 1.902045839092864
$ws.Range("C21").Value = 0.3827588639039732
$ws.Range("D21").Value = 0.1073500938741532
$ws.Range("E21").Value = 0.0564605383668404
$ws.Range("F21").Value = 2.601777283543157
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 1.861579153598399
$ws.Range("L21").Value = 0.241248838566392
$ws.Range("M21").Value = 0.3845065336347773
$ws.Range("B22").Value = 1.985546452005622
$ws.Range("C22").Value = 0.4107650891097023
$ws.Range("D22").Value = 0.1073933194174685
$ws.Range("E22").Value = 0.05682891407451773
$ws.Range("F22").Value = 2.659393261895957
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 1.887251037944282
$ws.Range("L22").Value = 0.2447180242865272
$ws.Range("M22").Value = 0.3978900228595847
$ws.Range("B23").Value = 1.940898085229605
$ws.Range("C23").Value = 0.3958126110863986
$ws.Range("D23").Value = 0.107367858401112
$ws.Range("E23").Value = 0.05663283616675585
$ws.Range("F23").Value = 2.628530884024883
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 1.87347876468921
$ws.Range("L23").Value = 0.2428568177431458
$ws.Range("M23").Value = 0.3907296590866807
$ws.Range("B24").Value = 1.773606619659063
$ws.Range("C24").Value = 0.3392888250024271
$ws.Range("D24").Value = 0.1073236192625018
$ws.Range("E24").Value = 0.05587846796509055
$ws.Range("F24").Value = 2.514081805330903
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 1.822862942176343
$ws.Range("L24").Value = 0.2360191033939429
$ws.Range("M24").Value = 0.3639909315497292
$ws.Range("B25").Value = 1.596981479964484
$ws.Range("C25").Value = 0.2785201283110439
$ws.Range("D25").Value = 0.1073853030541727
$ws.Range("E25").Value = 0.0550391094491296
$ws.Range("F25").Value = 2.395739331198484
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 1.771509277352607
$ws.Range("L25").Value = 0.2290924722067729
$ws.Range("M25").Value = 0.3359545370881989
